$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 36 (GARAMYCIN 0.1% CREAM 15 GM): balance ratio / selling price / transactions update ---
$ws.Range("H36").Value = "2:0"
$ws.Range("P36").Value = "66.0000"
$ws.Range("Q36").Value = "3:0"

# --- Row 70 (كريم فيبكس الازرق): balance ratio update ---
$ws.Range("H70").Value = "9:0"

# --- Insert two new data rows just above the totals row (currently row 73) ---
$ws.Rows("73:74").Insert()

# Copy formatting (styles/number formats/merges) from the row immediately above (row 72,
# still holding its original "مناديل جيب مبلله" formatting/pattern) down into the two
# freshly inserted blank rows 73:74.
$ws.Range("A72:Q72").Copy()
$ws.Range("A73:Q74").PasteSpecial(-4122)

# Re-create the per-row merged cell regions (lost by the format-only paste above).
$ws.Range("A73:B73").Merge()
$ws.Range("C73:G73").Merge()
$ws.Range("H73:K73").Merge()
$ws.Range("L73:M73").Merge()
$ws.Range("N73:O73").Merge()
$ws.Range("A74:B74").Merge()
$ws.Range("C74:G74").Merge()
$ws.Range("H74:K74").Merge()
$ws.Range("L74:M74").Merge()
$ws.Range("N74:O74").Merge()

# --- New row 71: مسواك اسنان (toothbrush) ---
$ws.Range("C71").Value = "مسواك اسنان "
$ws.Range("H71").Value = "4:0"
$ws.Range("L71").Value = "0"
$ws.Range("N71").Value = "15.00"
$ws.Range("P71").Value = "15.0000"
$ws.Range("Q71").Value = "1:0"

# --- New row 72: معجون سيجنال 25 مل (Signal toothpaste 25ml) ---
$ws.Range("C72").Value = "معجون سيجنال 25 مل"
$ws.Range("H72").Value = "6:0"
$ws.Range("L72").Value = "0"
$ws.Range("N72").Value = "20.00"
$ws.Range("P72").Value = "20.0000"
$ws.Range("Q72").Value = "1:0"

# --- Row 73 (newly added, carries forward the former row-71 content: معجون سيجنال اطفال 50 ملل) ---
$ws.Range("A73").Value = 67
$ws.Range("C73").Value = "معجون سيجنال اطفال 50 ملل"
$ws.Range("H73").Value = "2:0"
$ws.Range("L73").Value = "0"
$ws.Range("N73").Value = "55.00"
$ws.Range("P73").Value = "55.0000"
$ws.Range("Q73").Value = "1:0"

# --- Row 74 (newly added, carries forward the former row-72 content: مناديل جيب مبلله) ---
$ws.Range("A74").Value = 68
$ws.Range("C74").Value = "مناديل جيب مبلله "
$ws.Range("H74").Value = "9:0"
$ws.Range("L74").Value = "0"
$ws.Range("N74").Value = "6.00"
$ws.Range("P74").Value = "12.0000"
$ws.Range("Q74").Value = "2:0"

# --- Totals row (shifted from 73 to 75): update grand total ---
$ws.Range("P75").Value = 4404.7849999999999

# --- Footer row (shifted from 74 to 76): refresh the generated timestamp ---
$ws.Range("A76").Value = "Monday, 4 August, 2025 6:03 PM"
